$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = 65
$ws.Range("B66").Value = 1
$ws.Range("C66").Value = "2024-06-16 07:12:49"
$ws.Range("D66").Value = 200
$ws.Range("E66").Value = 8

$ws.Range("A67").Value = 66
$ws.Range("B67").Value = 2
$ws.Range("C67").Value = "2024-06-16 07:12:50"
$ws.Range("D67").Value = 200
$ws.Range("E67").Value = 0
